$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $text) {
    $cell = $ws.Range($ref)
    $cell.Value2 = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell "D2" '54.195.00'
Set-TextCell "E2" '  -8.02%  '
Set-TextCell "D3" '2.451.26'
Set-TextCell "E3" '  -13.26%  '
Set-TextCell "E4" '  +0.04%  '
Set-TextCell "D5" '465.64'
Set-TextCell "E5" '  -7.30%  '
Set-TextCell "D6" '131.83'
Set-TextCell "E6" '  -2.06%  '
Set-TextCell "D7" '0.995'
Set-TextCell "E7" '  -0.56%  '
Set-TextCell "D8" '0.489'
Set-TextCell "E8" '  -7.24%  '
Set-TextCell "D9" '2.471.97'
Set-TextCell "E9" '  -12.47%  '
Set-TextCell "D10" '0.0958'
Set-TextCell "E10" '  -6.22%  '
Set-TextCell "D11" '5.36'
Set-TextCell "E11" '  -9.11%  '
Set-TextCell "D12" '0.318'
Set-TextCell "E12" '  -8.05%  '
Set-TextCell "E13" '  -4.06%  '
Set-TextCell "D14" '2.883.19'
Set-TextCell "E14" '  -13.25%  '
Set-TextCell "D15" '54.420.99'
Set-TextCell "E15" '  -7.86%  '
Set-TextCell "E16" '  +1.62%  '
Set-TextCell "D17" '19.78'
Set-TextCell "E17" '  -7.53%  '
Set-TextCell "D18" '2.474.11'
Set-TextCell "E18" '  -12.46%  '
Set-TextCell "E19" '  -10.21%  '
Set-TextCell "D20" '312.01'
Set-TextCell "E20" '  -9.96%  '
Set-TextCell "D21" '9.44'
Set-TextCell "E21" '  -13.59%  '
Set-TextCell "E22" '  -0.59%  '
Set-TextCell "D23" '5.70'
Set-TextCell "E23" '  +1.16%  '
Set-TextCell "D24" '5.37'
Set-TextCell "E24" '  -13.42%  '
Set-TextCell "D25" '56.68'
Set-TextCell "E25" '  -9.90%  '
Set-TextCell "E26" '  +0.96%  '
Set-TextCell "B27" 'WrappedeETH'
Set-TextCell "C27" 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextCell "D27" '2.560.45'
Set-TextCell "E27" '  -13.14%  '
Set-TextCell "B28" 'Polygon'
Set-TextCell "C28" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell "D28" '0.384'
Set-TextCell "E28" '  -9.30%  '
Set-TextCell "E29" '  -8.26%  '
Set-TextCell "D30" '7.17'
Set-TextCell "E30" '  -1.72%  '
Set-TextCell "E31" '  -0.34%  '
Set-TextCell "D32" '0.0₃0730'
Set-TextCell "E32" '  -7.92%  '
Set-TextCell "D33" '152.08'
Set-TextCell "E33" '  +0.84%  '
Set-TextCell "D34" '17.71'
Set-TextCell "E34" '  -6.76%  '
Set-TextCell "D35" '1.43'
Set-TextCell "E35" '  -10.27%  '
Set-TextCell "D36" '5.01'
Set-TextCell "E36" '  -5.30%  '
Set-TextCell "E37" '  -13.40%  '
Set-TextCell "E38" '  -4.55%  '
Set-TextCell "D39" '0.793'
Set-TextCell "E39" '  -11.71%  '
Set-TextCell "D40" '33.59'
Set-TextCell "E40" '  -8.76%  '
Set-TextCell "B41" 'FirstDigitalUSD'
Set-TextCell "C41" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell "D41" '0.994'
Set-TextCell "E41" '  -0.50%  '
Set-TextCell "B42" 'Mantle'
Set-TextCell "C42" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell "D42" '0.605'
Set-TextCell "E42" '  -3.19%  '
Set-TextCell "D43" '0.0530'
Set-TextCell "E43" '  -4.25%  '
Set-TextCell "D44" '3.27'
Set-TextCell "E44" '  -6.27%  '
Set-TextCell "D45" '10.17'
Set-TextCell "E45" '  -1.74%  '
Set-TextCell "D46" '1.24'
Set-TextCell "E46" '  -6.97%  '
Set-TextCell "D47" '1.965.68'
Set-TextCell "E47" '  -11.81%  '
Set-TextCell "E48" '  -0.75%  '
Set-TextCell "D49" '0.0870'
Set-TextCell "E49" '  -1.60%  '
Set-TextCell "D50" '4.34'
Set-TextCell "E50" '  -4.51%  '
Set-TextCell "D51" '16.71'
Set-TextCell "E51" '  -12.41%  '
